$d = $word.ActiveDocument

# 1) Remove the empty paragraph that only contained the "_GoBack" bookmark,
#    leaving a plain empty paragraph in its place.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Remove the paragraph advertising the MRIQC python gist script, together
#    with the blank paragraph that followed it (the paragraph ends up
#    directly followed by the "Now you're ready for any BIDS app ;-)" text,
#    matching the already-existing blank line before it).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*python script that runs MRIQC on your BIDS dataset*") {
        $target = $i
        break
    }
}
if ($target -ne $null) {
    $pStart = $d.Paragraphs.Item($target)
    $pNext = $d.Paragraphs.Item($target + 1)
    $r = $d.Range($pStart.Range.Start, $pNext.Range.End)
    $r.Delete()
}

# 3) Point the "Credit to many of the scripts found here" hyperlink at the
#    same target used by the (now removed) gist link, instead of the
#    bids-starter-kit matlabCode tree, while keeping its displayed text.
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.Address -eq "https://github.com/bids-standard/bids-starter-kit/tree/master/matlabCode") {
        $h.Address = "https://gist.github.com/marcoaqil/c0e0584513fd482d6ea8e9b164b7c1f4"
        break
    }
}
